# Update TPM-derived NATMI ligand-receptor metrics on Sheet1
# (rows 2-10, columns G/H/I/J/M/N/O/P/Q/R/S/T) with refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 35.906979
$ws.Range("H2").Value = 107.720937
$ws.Range("I2").Value = 0.6107087147789413
$ws.Range("J2").Value = 0.6107087147789412
$ws.Range("M2").Value = 1.116695
$ws.Range("N2").Value = 3.350085
$ws.Range("O2").Value = 0.008174214292497491
$ws.Range("P2").Value = 0.008174214292497492
$ws.Range("Q2").Value = 40.097143914405
$ws.Range("R2").Value = 360.874295229645
$ws.Range("S2").Value = 0.004992063904898795
$ws.Range("T2").Value = 0.004992063904898795
$ws.Range("G3").Value = 35.906979
$ws.Range("H3").Value = 107.720937
$ws.Range("I3").Value = 0.6107087147789413
$ws.Range("J3").Value = 0.6107087147789412
$ws.Range("O3").Value = 0.8193429796700005
$ws.Range("P3").Value = 0.8193429796700005
$ws.Range("Q3").Value = 4019.14021280786
$ws.Range("R3").Value = 36172.26191527073
$ws.Range("S3").Value = 0.5003798980774142
$ws.Range("T3").Value = 0.5003798980774141
$ws.Range("G4").Value = 35.906979
$ws.Range("H4").Value = 107.720937
$ws.Range("I4").Value = 0.6107087147789413
$ws.Range("J4").Value = 0.6107087147789412
$ws.Range("O4").Value = 0.172482806037502
$ws.Range("P4").Value = 0.1724828060375021
$ws.Range("Q4").Value = 846.083507107695
$ws.Range("R4").Value = 7614.751563969255
$ws.Range("S4").Value = 0.1053367527966283
$ws.Range("T4").Value = 0.1053367527966283
$ws.Range("I5").Value = 0.2899643113254147
$ws.Range("J5").Value = 0.2899643113254147
$ws.Range("M5").Value = 1.116695
$ws.Range("N5").Value = 3.350085
$ws.Range("O5").Value = 0.008174214292497491
$ws.Range("P5").Value = 0.008174214292497492
$ws.Range("Q5").Value = 19.03811168875334
$ws.Range("R5").Value = 171.34300519878
$ws.Range("S5").Value = 0.002370230417950397
$ws.Range("T5").Value = 0.002370230417950397
$ws.Range("I6").Value = 0.2899643113254147
$ws.Range("J6").Value = 0.2899643113254147
$ws.Range("O6").Value = 0.8193429796700005
$ws.Range("P6").Value = 0.8193429796700005
$ws.Range("S6").Value = 0.237580222839325
$ws.Range("T6").Value = 0.2375802228393249
$ws.Range("I7").Value = 0.2899643113254147
$ws.Range("J7").Value = 0.2899643113254147
$ws.Range("O7").Value = 0.172482806037502
$ws.Range("P7").Value = 0.1724828060375021
$ws.Range("S7").Value = 0.05001385806813936
$ws.Range("T7").Value = 0.05001385806813936
$ws.Range("I8").Value = 0.09932697389564409
$ws.Range("J8").Value = 0.09932697389564407
$ws.Range("M8").Value = 1.116695
$ws.Range("N8").Value = 3.350085
$ws.Range("O8").Value = 0.008174214292497491
$ws.Range("P8").Value = 0.008174214292497492
$ws.Range("Q8").Value = 6.521485399660001
$ws.Range("R8").Value = 58.69336859694
$ws.Range("S8").Value = 0.0008119199696482991
$ws.Range("T8").Value = 0.0008119199696482991
$ws.Range("I9").Value = 0.09932697389564409
$ws.Range("J9").Value = 0.09932697389564407
$ws.Range("O9").Value = 0.8193429796700005
$ws.Range("P9").Value = 0.8193429796700005
$ws.Range("S9").Value = 0.08138285875326139
$ws.Range("T9").Value = 0.08138285875326137
$ws.Range("I10").Value = 0.09932697389564409
$ws.Range("J10").Value = 0.09932697389564407
$ws.Range("O10").Value = 0.172482806037502
$ws.Range("P10").Value = 0.1724828060375021
$ws.Range("S10").Value = 0.01713219517273441
$ws.Range("T10").Value = 0.01713219517273441
